$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-06 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-07 Thursday", 2) | Out-Null
$d.Content.Find.Execute("83×21=", $true, $false, $false, $false, $false, $true, 1, $false, "37×55=", 2) | Out-Null
$d.Content.Find.Execute("89×94=", $true, $false, $false, $false, $false, $true, 1, $false, "82×70=", 2) | Out-Null
$d.Content.Find.Execute("61×37=", $true, $false, $false, $false, $false, $true, 1, $false, "23×12=", 2) | Out-Null
$d.Content.Find.Execute("26×67=", $true, $false, $false, $false, $false, $true, 1, $false, "82×30=", 2) | Out-Null
$d.Content.Find.Execute("64×36=", $true, $false, $false, $false, $false, $true, 1, $false, "20×72=", 2) | Out-Null
$d.Content.Find.Execute("65×74=", $true, $false, $false, $false, $false, $true, 1, $false, "81×56=", 2) | Out-Null
$d.Content.Find.Execute("71×77=", $true, $false, $false, $false, $false, $true, 1, $false, "84×76=", 2) | Out-Null
$d.Content.Find.Execute("56×69=", $true, $false, $false, $false, $false, $true, 1, $false, "38×53=", 2) | Out-Null
$d.Content.Find.Execute("81×88=", $true, $false, $false, $false, $false, $true, 1, $false, "97×16=", 2) | Out-Null
$d.Content.Find.Execute("63×36=", $true, $false, $false, $false, $false, $true, 1, $false, "66×14=", 2) | Out-Null
$d.Content.Find.Execute("45×81=", $true, $false, $false, $false, $false, $true, 1, $false, "20×75=", 2) | Out-Null
$d.Content.Find.Execute("53×46=", $true, $false, $false, $false, $false, $true, 1, $false, "77×77=", 2) | Out-Null
$d.Content.Find.Execute("94×76=", $true, $false, $false, $false, $false, $true, 1, $false, "90×32=", 2) | Out-Null
$d.Content.Find.Execute("22×41=", $true, $false, $false, $false, $false, $true, 1, $false, "24×25=", 2) | Out-Null
$d.Content.Find.Execute("38×72=", $true, $false, $false, $false, $false, $true, 1, $false, "20×79=", 2) | Out-Null
$d.Content.Find.Execute("17×24=", $true, $false, $false, $false, $false, $true, 1, $false, "84×44=", 2) | Out-Null
$d.Content.Find.Execute("53×34=", $true, $false, $false, $false, $false, $true, 1, $false, "67×19=", 2) | Out-Null
$d.Content.Find.Execute("20×78=", $true, $false, $false, $false, $false, $true, 1, $false, "79×88=", 2) | Out-Null
$d.Content.Find.Execute("76×57=", $true, $false, $false, $false, $false, $true, 1, $false, "70×91=", 2) | Out-Null
$d.Content.Find.Execute("33×65=", $true, $false, $false, $false, $false, $true, 1, $false, "28×79=", 2) | Out-Null
$d.Content.Find.Execute("34×84=", $true, $false, $false, $false, $false, $true, 1, $false, "72×40=", 2) | Out-Null
$d.Content.Find.Execute("57×27=", $true, $false, $false, $false, $false, $true, 1, $false, "84×82=", 2) | Out-Null
$d.Content.Find.Execute("31×29=", $true, $false, $false, $false, $false, $true, 1, $false, "87×34=", 2) | Out-Null
$d.Content.Find.Execute("11×85=", $true, $false, $false, $false, $false, $true, 1, $false, "98×93=", 2) | Out-Null
$d.Content.Find.Execute("39×76=", $true, $false, $false, $false, $false, $true, 1, $false, "23×54=", 2) | Out-Null
